# Double every value in column B (rows 2 through 463), which holds the
# "Average Coefficient Of Friction B" series, leaving column A (time) and
# the header row untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.Value2 = $cell.Value2 * 2
}
